# Eliminated some compiler warning messages.
#
# The source data picked up one more counted item (I10 = 2), which in turn
# changes the rolled-up totals / derived rows below it (I39:I42). The sheet
# view's scroll position and active selection are also updated to reflect
# where the author was working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data point - number in the "I" column (e.g. estimated/actual hours)
# for row 10 (IST package). This feeds SUM(I2:I37) in I39 and the derived
# rows I40 (I39/5), I41 (I40*7/5) and I42 (I41/30), which Excel will
# recalculate automatically.
$ws.Range("I10").Value = 2

# Make sure everything is recalculated before saving.
$excel.CalculateFull()

# Update the active selection to match where the author ended up.
$ws.Range("J16").Select()
